$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should look exactly like the
# existing header cells (bold, bordered, centered style used by B1:H1).
# Copy the formatting from H1 so the new header cells share the same style
# index instead of minting a new one, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-10.
$values = @{
    2  = @(6, 6)
    3  = @(8, 8)
    4  = @(9, 9)
    5  = @(7, 7)
    6  = @(8, 8)
    7  = @(7, 8)
    8  = @(7, 7)
    9  = @(9, 9)
    10 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
